$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "Protein"
$ws.Range("I1").Value = "Ontogeny"
$ws.Range("I2").Select() | Out-Null
